# Business configuration.xlsx - "Latest commit after demo"
# Update the SetUp sheet credentials, remove the now-unused extra data row,
# and leave the workbook focused on the SetUp sheet (cell A3 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SetUp")

# New login credentials replace the old ones
$ws.Range("A2").Value = "Vikesh.patil@aimdek.com"
$ws.Range("B2").Value = "Vikesh@1989"

# The stray numeric row (A3:B3 = 112 / 336) is no longer needed
$ws.Range("A3:B3").Clear()

# Make "SetUp" the active sheet/tab, with A3 as the selected cell
$ws.Activate()
$ws.Range("A3").Select() | Out-Null
